# Update the 'F' column (attendance/view counts) values across all four
# worksheets of the workbook, matching the regenerated gh-pages data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 609
$ws.Cells.Item(5, 6).Value = 2692
$ws.Cells.Item(9, 6).Value = 262
$ws.Cells.Item(10, 6).Value = 6009
$ws.Cells.Item(11, 6).Value = 65
$ws.Cells.Item(13, 6).Value = 4921
$ws.Cells.Item(15, 6).Value = 93
$ws.Cells.Item(16, 6).Value = 6
$ws.Cells.Item(17, 6).Value = 2527
$ws.Cells.Item(18, 6).Value = 1318
$ws.Cells.Item(19, 6).Value = 492
$ws.Cells.Item(20, 6).Value = 1201
$ws.Cells.Item(22, 6).Value = 275
$ws.Cells.Item(23, 6).Value = 110
$ws.Cells.Item(26, 6).Value = 218
$ws.Cells.Item(29, 6).Value = 1337
$ws.Cells.Item(32, 6).Value = 2070
$ws.Cells.Item(33, 6).Value = 281
$ws.Cells.Item(34, 6).Value = 562
$ws.Cells.Item(35, 6).Value = 54
$ws.Cells.Item(36, 6).Value = 235
$ws.Cells.Item(37, 6).Value = 1457
$ws.Cells.Item(38, 6).Value = 606
$ws.Cells.Item(40, 6).Value = 545
$ws.Cells.Item(41, 6).Value = 240
$ws.Cells.Item(42, 6).Value = 1712
$ws.Cells.Item(43, 6).Value = 2501
$ws.Cells.Item(45, 6).Value = 104
$ws.Cells.Item(47, 6).Value = 94
$ws.Cells.Item(48, 6).Value = 65
$ws.Cells.Item(49, 6).Value = 89

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 16
$ws.Cells.Item(8, 6).Value = 309
$ws.Cells.Item(9, 6).Value = 166
$ws.Cells.Item(10, 6).Value = 77
$ws.Cells.Item(11, 6).Value = 197
$ws.Cells.Item(19, 6).Value = 35
$ws.Cells.Item(22, 6).Value = 318
$ws.Cells.Item(23, 6).Value = 326

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(6, 6).Value = 1688
$ws.Cells.Item(7, 6).Value = 562
$ws.Cells.Item(8, 6).Value = 1431
$ws.Cells.Item(9, 6).Value = 1798
$ws.Cells.Item(10, 6).Value = 2386
$ws.Cells.Item(11, 6).Value = 793
$ws.Cells.Item(12, 6).Value = 680

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 1688
$ws.Cells.Item(5, 6).Value = 609
$ws.Cells.Item(6, 6).Value = 562
$ws.Cells.Item(7, 6).Value = 2692
$ws.Cells.Item(9, 6).Value = 1431
$ws.Cells.Item(10, 6).Value = 2386
$ws.Cells.Item(11, 6).Value = 6010
$ws.Cells.Item(12, 6).Value = 793
$ws.Cells.Item(14, 6).Value = 16
$ws.Cells.Item(15, 6).Value = 65
$ws.Cells.Item(16, 6).Value = 4921
$ws.Cells.Item(17, 6).Value = 93
$ws.Cells.Item(18, 6).Value = 2527
$ws.Cells.Item(19, 6).Value = 1318
$ws.Cells.Item(20, 6).Value = 492
$ws.Cells.Item(21, 6).Value = 1201
$ws.Cells.Item(22, 6).Value = 275
$ws.Cells.Item(23, 6).Value = 110
$ws.Cells.Item(25, 6).Value = 166
$ws.Cells.Item(26, 6).Value = 218
$ws.Cells.Item(28, 6).Value = 1337
$ws.Cells.Item(29, 6).Value = 2070
$ws.Cells.Item(30, 6).Value = 281
$ws.Cells.Item(31, 6).Value = 562
$ws.Cells.Item(32, 6).Value = 235
$ws.Cells.Item(34, 6).Value = 1457
$ws.Cells.Item(35, 6).Value = 606
$ws.Cells.Item(37, 6).Value = 545
$ws.Cells.Item(39, 6).Value = 318
$ws.Cells.Item(40, 6).Value = 240
$ws.Cells.Item(42, 6).Value = 1712
$ws.Cells.Item(43, 6).Value = 2501
$ws.Cells.Item(44, 6).Value = 104
$ws.Cells.Item(46, 6).Value = 94
$ws.Cells.Item(47, 6).Value = 65

